# Chat Relay - Design Presentation.pptx
#
# "Design Presentation - Fixed spelling issues and moving around slides"
#
# The only audience-visible structural change in this revision is a swap of
# two adjacent slides: "Chat / Message  Descriptions" (previously slide 6)
# now comes *before* "User / IT Admin Comparison" (previously slide 5).
#
# (The deck's cached `datetimeFigureOut` footer text on the Handout Master /
# Notes Master - "4/9/2025" -> "4/10/25" - is a PowerPoint auto-refresh of an
# auto date field driven by the system clock at save time, not a deliberate
# edit, and the field glyphs themselves aren't user-editable text.)

$p = $ppt.ActivePresentation

# Move slide 6 ("Chat / Message  Descriptions") up to position 5; this
# pushes slide 5 ("User / IT Admin Comparison") down to position 6 -
# the same adjacent-pair swap shown in the diff's slide-order list.
$p.Slides.Item(6).MoveTo(5)
